$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "java -jar .\baron_paul_client.jar add [google|microsoft] <user>"
#    -> whole line bold, append " <userKey>"
# ---------------------------------------------------------------------
$pAdd = $d.Paragraphs(6)
$pAdd.Range.Find.Execute("<user>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<user> <userKey>", 2)
$pAdd.Range.Bold = 1

# ---------------------------------------------------------------------
# 2) "A FAIRE : Permet d'ajouter un compte Google / Microsoft à l'<user>"
#    -> drop the "A FAIRE : " prefix, <user> -> <userKey>
# ---------------------------------------------------------------------
$pPermet = $d.Paragraphs(7)
$pPermet.Range.Find.Execute("A FAIRE : ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "", 2)
$pPermet.Range.Find.Execute("<user>", $false, $false, $false, $false, `
    $false, $true, 1, $false, "<userKey>", 2)

# ---------------------------------------------------------------------
# 3) "java -jar .\baron_paul_client.jar view <user>"
#    -> whole line bold, <user> -> <userKey>
# ---------------------------------------------------------------------
$pView = $d.Paragraphs(8)
$pView.Range.Find.Execute("<user>", $false, $false, $false, $false, $false, `
    $true, 1, $false, "<userKey>", 2)
$pView.Range.Bold = 1

# ---------------------------------------------------------------------
# 4) Append " Google + Microsoft" to the three result bullet lines
# ---------------------------------------------------------------------
foreach ($idx in 10, 11, 12) {
    $p = $d.Paragraphs($idx)
    $r = $p.Range
    $endRange = $d.Range($r.End - 1, $r.End - 1)
    $endRange.InsertAfter(" Google + Microsoft")
}

# ---------------------------------------------------------------------
# 5) "java -jar .\baron_paul_client.jar email <user>" -> bold whole line
# ---------------------------------------------------------------------
$d.Paragraphs(13).Range.Bold = 1

# ---------------------------------------------------------------------
# 6) "java -jar .\baron_paul_client.jar contact <user>" -> bold whole line
# ---------------------------------------------------------------------
$d.Paragraphs(15).Range.Bold = 1

# ---------------------------------------------------------------------
# 7) "java -jar .\baron_paul_client.jar calendar <user>" -> bold whole line
# ---------------------------------------------------------------------
$d.Paragraphs(17).Range.Bold = 1

# ---------------------------------------------------------------------
# 8) New warning paragraph at the very end of the document (before the
#    trailing blank paragraph's sectPr), bold "/ ! \" bookends around a
#    plain-text warning about token expiry.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$nr = $newPara.Range
$fullText = "/ ! \ Les tokens des comptes ajoutés expirent rapidement ! Pensez à les renouveler en cas de problème ; supprimez les comptes Microsoft dans la base de données et supprimez le ‘StoredCredential’ dans votre répertoire Google défini. / ! \"
$nr.Text = $fullText

$pStart = $newPara.Range.Start
$boldStart = $d.Range($pStart, $pStart + 5)
$boldStart.Bold = 1

$pEnd = $newPara.Range.End - 1
$boldEnd = $d.Range($pEnd - 5, $pEnd)
$boldEnd.Bold = 1

Write-Output "done"
